$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextCell "D2" "63.613.22"
Set-TextCell "E2" "  +4.66%  "

# Row 3 - Ethereum
Set-TextCell "D3" "2.478.79"
Set-TextCell "E3" "  +5.90%  "

# Row 4 - TetherUSD
Set-TextCell "E4" "  +0.15%  "

# Row 5 - BNB
Set-TextCell "D5" "568.89"
Set-TextCell "E5" "  +3.73%  "

# Row 6 - Solana
Set-TextCell "D6" "143.60"
Set-TextCell "E6" "  +9.07%  "

# Row 7 - USDC
Set-TextCell "E7" "  +0.11%  "

# Row 8 - XRP
Set-TextCell "D8" "0.592"
Set-TextCell "E8" "  +2.38%  "

# Row 9 - LidoStakedEther
Set-TextCell "D9" "2.478.01"
Set-TextCell "E9" "  +5.97%  "

# Row 10 - Dogecoin
Set-TextCell "E10" "  +4.46%  "

# Row 11 - Toncoin
Set-TextCell "E11" "  +4.40%  "

# Row 12 - TRON
Set-TextCell "E12" "  +1.09%  "

# Row 13 - Cardano
Set-TextCell "E13" "  +5.12%  "

# Row 14 - Avalanche
Set-TextCell "D14" "26.46"
Set-TextCell "E14" "  +12.07%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextCell "D15" "2.924.93"
Set-TextCell "E15" "  +6.11%  "

# Row 16 - WrappedBTC
Set-TextCell "D16" "63.494.37"
Set-TextCell "E16" "  +4.59%  "

# Row 17 - ShibaInu
Set-TextCell "E17" "  +6.50%  "

# Row 18 - WrappedEther
Set-TextCell "D18" "2.481.33"
Set-TextCell "E18" "  +6.18%  "

# Row 19 - Chainlink
Set-TextCell "D19" "11.32"
Set-TextCell "E19" "  +6.26%  "

# Row 20 - BitcoinCash
Set-TextCell "D20" "342.58"
Set-TextCell "E20" "  +8.67%  "

# Row 21 - Polkadot
Set-TextCell "D21" "4.32"
Set-TextCell "E21" "  +5.46%  "

# Row 22 - Uniswap
Set-TextCell "E22" "  +4.19%  "

# Row 23 - Dai
Set-TextCell "E23" "  -0.01%  "

# Row 24 - Litecoin
Set-TextCell "D24" "65.83"
Set-TextCell "E24" "  +3.00%  "

# Row 25 - Kaspa
Set-TextCell "E25" "  +1.38%  "

# Row 26 - Binance-PegBSC-USD
Set-TextCell "E26" "  +0.04%  "

# Row 27 - Fetch.AI
Set-TextCell "E27" "  +8.74%  "

# Row 28 - InternetComputer(DFINITY)
Set-TextCell "D28" "8.24"
Set-TextCell "E28" "  +3.83%  "

# Row 29 - SuiNetwork
Set-TextCell "E29" "  +5.99%  "

# Row 30 - PEPE
Set-TextCell "D30" "0.0₃0828"
Set-TextCell "E30" "  +12.61%  "

# Row 31 - Aptos
Set-TextCell "D31" "6.87"
Set-TextCell "E31" "  +15.05%  "

# Row 32 - PancakeSwap
Set-TextCell "E32" "  +6.42%  "

# Row 33 - Monero
Set-TextCell "D33" "177.32"
Set-TextCell "E33" "  +2.59%  "

# Row 34 - ImmutableX
Set-TextCell "E34" "  +10.74%  "

# Row 35 - PolygonEcosystemToken
Set-TextCell "E35" "  +4.21%  "

# Row 36 - EthereumClassic
Set-TextCell "D36" "19.01"
Set-TextCell "E36" "  +5.65%  "

# Row 37 - Bittensor
Set-TextCell "D37" "372.27"
Set-TextCell "E37" "  +13.38%  "

# Row 38 - NEARProtocol
Set-TextCell "D38" "4.48"
Set-TextCell "E38" "  +8.13%  "

# Row 39 - USDe (unchanged)

# Rows 40 & 41 swap: Stacks <-> FirstDigitalUSD
Set-TextCell "B40" "FirstDigitalUSD"
Set-TextCell "C40" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCell "D40" "1.00"
Set-TextCell "E40" "  +0.17%  "

Set-TextCell "B41" "Stacks"
Set-TextCell "C41" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell "D41" "1.72"
Set-TextCell "E41" "  +11.52%  "

# Row 42 - OKB
Set-TextCell "D42" "40.46"
Set-TextCell "E42" "  +6.02%  "

# Row 43 - Aave
Set-TextCell "D43" "151.47"
Set-TextCell "E43" "  +10.45%  "

# Row 44 - Filecoin
Set-TextCell "D44" "3.73"
Set-TextCell "E44" "  +6.48%  "

# Row 45 - InjectiveProtocol
Set-TextCell "E45" "  +8.61%  "

# Row 46 - Mantle
Set-TextCell "E46" "  +5.69%  "

# Row 47 - Stellar
Set-TextCell "D47" "0.0967"
Set-TextCell "E47" "  +2.56%  "

# Row 48 - Hedera
Set-TextCell "D48" "0.0525"
Set-TextCell "E48" "  +5.46%  "

# Row 49 - BabyDogeCoin
Set-TextCell "D49" "0.0₆0239"
Set-TextCell "E49" "  +8.46%  "

# Row 50 - VeChain
Set-TextCell "D50" "0.0227"
Set-TextCell "E50" "  +4.48%  "

# Row 51 - EnergySwap
Set-TextCell "E51" "  +7.24%  "
